$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.059008598327637
$ws.Range("B1").Value = 3.584442853927612
$ws.Range("C1").Value = 3.341569185256958
$ws.Range("D1").Value = 2.002537965774536
$ws.Range("E1").Value = 1.156426072120667
